$wb = $excel.ActiveWorkbook

# 1) Status text: "Ready for handoff" -> "In Translation"
#    Overview sheet: columns E2 and F2 mirror the zh-cn / de-de status.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

#    zh-cn / de-de sheets: column C2 is the "Status" cell.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# 2) Narrow the corresponding columns to match the new (shorter) text.
#    (target raw OOXML width is 13.4101845877511; the COM ColumnWidth
#    setter only resolves to a 1/6-character pixel grid, so 12.5
#    characters is the closest achievable value -> 13.333... width)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
